$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refresh of the mesh / device inventory table (commit: "add mesh RE script").
# Columns: A=flag(0), B=DeviceID, C=SN, D=MAC, E=Firmware, F=Model, G=Device-Type
$data = @(
    @(0, "QUwSDcl2ASZFQANSoKgOjUGUXQKpxwjzP1kVJAXhvzHLrhEPt6dw4C0Tha_3dm7W", "2219111B02429", "E4C32AD28FC4", "2.5.14 Build 20240827 Rel. 31177", "Deco M5", "DECO"),
    @(0, "z3UV9QwRYB0a-cW7zheCm5-w3OYecpCkiOOSvRPsjmXVVNlW8ieW2bbXAaj7v8on", "2237193A00760", "5091E3E64470", "2.0.0 Build 20220408 Rel. 70690", "Deco X50", "DECO"),
    @(0, "-pSp3cJWLo5jIJucwoWgqYfFzu2tJQcXDruTMCfFZP3KbMObQdZQnnWwOwSiYuQA", "22221K3A00172", "54AF97557F2C", "2.0.0 Build 20220408 Rel. 70690", "Deco X50", "DECO"),
    @(0, "u0XTyvoaV9Hz4_FLE3Rr8sqrJD6cDfbRFkX226tBpM-j0rS91yfPnmA_Volthc3v", "22465B9000110", "F0090D676B60", "0.1.0 3.0.0 v60d5.0 Build 240918 Rel.44802n_Beta", "HB210 Pro", "AGINET"),
    @(0, "j5zFq07CKJXQM80HePqC5RK80L-gh1rJh6IiMKGcClO-DfJVpEgMMLbKU0FXx7pu", "2241584000038", "74FECEDAE1E2", "0.1.0 3.0.0 v60af.0 Build 240424 Rel.48510n", "HB610V2", "AGINET"),
    @(0, "CwkNE9l2s_NOUoUxYFLjai_bEQyJG1nZ0W-UgcBMKcEUrS3k3Ql7QwBvJnN1ugEs", "2245034R00016", "6083E7C01398", "0.2.0 3.0.0 v60be.0 Build 240906 Rel.40026n", "HB710", "AGINET"),
    @(0, "vZC-Qg9GyIDl32SE44YddRmAP6_nKg0Lf1O4Fh_Yo8Pzu-TUQrUAD4DhMOXTCbdB", "2243627000112", "2023515A2B98", "0.2.0 3.0.0 v60be.0 Build 240906 Rel.40026n", "HB710", "AGINET"),
    @(0, "CjD6AafkLUhwoLzqgIW8JKgewZwGF6VI4bW2sI4A0qs7g1ZGuXO6rbbCNpel0f0S", "22482S0000148", "7CF17E38C6EE", "0.4.0 3.0.0 v6095.0 Build 240812 Rel.12929n", "HB810", "AGINET"),
    @(0, "piuNbmyRSm29ORfCSk9ROSQIb5CAJsh0pOdrcK4-wl0fllzE1ws4xBDNaSzeYoW0", "22346B4000012", "5C628BF46D0F", "0.3.0 3.0.0 v6095.0 Build 240509 Rel.21531n", "HB810", "AGINET"),
    @(0, "h07Ab5DSHFbSpJMZhJA593mBbBUmw5Rl5Au_FWK4FsfzGWIQ2cFNXxNGI_2Fcmak", "223A4B1000018", "A842A1A4C394", "0.3.0 3.0.0 v6095.0 Build 240509 Rel.21531n", "HB810", "AGINET"),
    @(0, "TbUKc1W_jAJfJC5Yb74-Edj2N_yRESFXKqcc0QQGHNko0nfzRdu5HjxiZe6Qz0Wd", "22360N3001322", "40ED00CB1D2F", "0.16.0 2.0.0 v605f.0 Build 240430 Rel.29144n", "HX220", "AGINET"),
    @(0, "TCD7FsU6vv5SwER9OVg3ZEg4PuuZVWRFpRgz5Wz6NL6X9H_q4Xx0RYO_2B3OnKUc", "22360N3001039", "40ED00CB19DE", "0.16.0 2.0.0 v605f.0 Build 240430 Rel.29144n", "HX220", "AGINET"),
    @(0, "mPdUSXiIMTYJMIQJeJRCKZMsVwLz54pnMuCzZyfDpotOMkSamXnXR7IO1MdQ2WmX", "22431Q4000183", "40AE30999B7A", "0.9.0 3.0.0 v6065.0 Build 240806 Rel.46875n", "HX510", "AGINET"),
    @(0, "9yzsY0Vi4oNGUEZNBRnwX3erSH5L9AllmNC_3MgxmSWQiloBzq8wANdrrjFcoNon", "22360N3001370", "40ED00CB1DBF", "0.17.0 2.0.0 v605f.0 Build 240522 Rel.4911n", "HX220", "AGINET"),
    @(0, "JmcIw16f6q6eF_PMVmVk-WHN5DvNeKo4c3Fw0YKqBm6CqHqnXwYBBva3lSt8C-0y", "223CACS000016", "20362690B8F5", "0.1.0 3.1.9 v60b3.0 Build 240515 Rel.78630n", "EX921v", "AGINET"),
    @(0, "mn4qqGRXlEsgxVxPuZ-wFElcJiV0_wiHpkfrsWP1Mx-ccFel97uAbyBZuuPIEpI7", "22334M3002425", "788CB5C7D660", "0.6.0 3.0.0 v6066.0 Build 220630 Rel.56119n", "XX230v", "AGINET"),
    @(0, "K4ouszgY_l92YuT_5yqdOiHcZzfBFzmBZEtZFs5RXrI1dQsJnIJggleK38pHp_Z6", "223A4B1000236", "A842A1A4CC18", "0.3.0 3.0.0 v6095.0 Build 240509 Rel.21531n", "HB810", "AGINET"),
    @(0, "VtSm-VcF7DindwuZctOpTjtEA0gbkphdFO3LBVATPJu54RydB9udoWrgwupybwh2", "2245034R00023", "6083E7C013D0", "0.2.0 3.0.0 v60be.0 Build 240906 Rel.40026n", "HB710", "AGINET"),
    @(0, "oxn8uthuVz7ZhR1mI-dN0_DPoZd1hqNrZrJ1E3o5M3gSEsZTMWSXNxazgN7tkhoC", "2241584000037", "74FECEDAE1D8", "0.1.0 3.0.0 v60af.0 Build 240424 Rel.48510n", "HB610V2", "AGINET"),
    @(0, "dL_Feznnvl_Eq6SiQ8nfkiujXj3GFfZGwmRxByPcUT6EOkxTqFW4T_J9rB6Si_bP", "2241542000290", "74FECECE54CB", $null, "HX510-PoE", "AGINET")
)

# Scratch cell, far outside the used range, that still carries the sheet's
# untouched default style (numFmtId 0 / no explicit xf). Used to "reset"
# the format of cells whose text would otherwise be auto-detected as a
# number by Excel, without creating any new style entries.
$blank = $ws.Cells.Item(500, 500)

function Set-TextValue($cell, [string]$text) {
    $looksNumeric = $text -match '^[0-9]+$'
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value2 = $text
        $blank.Copy() | Out-Null
        $cell.PasteSpecial(-4122) | Out-Null
        $excel.CutCopyMode = $false
    } else {
        $cell.Value2 = $text
    }
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rec = $data[$i]

    # Rows beyond the original extent (20, 21) don't carry the bordered/bold
    # style on column A yet -- clone it from the row above before writing.
    if ($row -gt 19) {
        $ws.Cells.Item(19, 1).Copy() | Out-Null
        $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
        $excel.CutCopyMode = $false
    }

    $ws.Cells.Item($row, 1).Value2 = $rec[0]
    Set-TextValue $ws.Cells.Item($row, 2) $rec[1]
    Set-TextValue $ws.Cells.Item($row, 3) $rec[2]
    Set-TextValue $ws.Cells.Item($row, 4) $rec[3]
    if ($rec[4] -ne $null) {
        Set-TextValue $ws.Cells.Item($row, 5) $rec[4]
    }
    Set-TextValue $ws.Cells.Item($row, 6) $rec[5]
    Set-TextValue $ws.Cells.Item($row, 7) $rec[6]
}
